$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.203.13'
$ws.Range('E2').Value = '  -1.93%  '
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''306.62'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''0.5189'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.73%  '
$ws.Range('D8').Value = '''0.3739'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').Value = '''0.07160'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').Value = '''0.8937'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D11').Value = '''20.77'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.24%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.871.23'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.07535'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.45%  '
$ws.Range('E14').Value = '  -2.55%  '
$ws.Range('D15').Value = '''90.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.89%  '
$ws.Range('D16').Value = '''1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '''0.000008510'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('E18').Value = '  -2.47%  '
$ws.Range('D19').Value = '''0.9996'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '27.231.23'
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').Value = '''5.010'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.50%  '
$ws.Range('D22').Value = '2.107.85'
$ws.Range('E22').Value = '  -2.69%  '
$ws.Range('D23').Value = '''10.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.62%  '
$ws.Range('D24').Value = '''6.477'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('D25').Value = '''1.836'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.64%  '
$ws.Range('D26').Value = '''145.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.84%  '
$ws.Range('D27').Value = '''18.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.58%  '
$ws.Range('D28').Value = '''2.089'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.51%  '
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('D30').Value = '''4.664'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('D31').Value = '''4.685'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.15%  '
$ws.Range('D32').Value = '''0.09264'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.02%  '
$ws.Range('D33').Value = '''0.05138'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.87%  '
$ws.Range('D34').Value = '''3.085'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('D35').Value = '''1.161'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.92%  '
$ws.Range('D36').Value = '''0.7280'
$ws.Range('D36').Style = 'Normal'
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '''3.129'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.91%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.02032'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.98%  '
$ws.Range('D39').Value = '''2.516'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').Value = '''0.5318'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.13%  '
$ws.Range('D42').Value = '''6.540'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.90%  '
$ws.Range('D43').Value = '''116.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').Value = '''8.347'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.91%  '
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('D46').Value = '''0.4635'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.73%  '
$ws.Range('D47').Value = '''0.9998'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('E48').Value = '  -4.06%  '
$ws.Range('E49').Value = '  -2.85%  '
$ws.Range('D50').Value = '''36.81'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').Value = '''63.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.52%  '
